$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1), columns A-F
$ws.Range("A1").Value = "First name"
$ws.Range("B1").Value = "Last name"
$ws.Range("C1").Value = "Bank name"
$ws.Range("D1").Value = "Bank account no."
$ws.Range("E1").Value = "Amount"
$ws.Range("F1").Value = "Fundraiser ID"

# Data (row 2), columns A-F - text-like numeric values stored as text
$ws.Range("A2").Value = "Fasikaw"
$ws.Range("B2").Value = "Kindye"
$ws.Range("C2").Value = "Cvg"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "123456789"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "10.00"
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").Value = "6155f2754bde6bb71afe5f7a"

# Newly added "status" column (G) - appended after existing data
$ws.Range("G1").Value = "status"
$ws.Range("G2").Value = 1

# Column width for column F (26 "characters" once persisted)
$ws.Columns.Item(6).ColumnWidth = 25.16666666666667

$ws.Range("G2").Select()
